$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "ACORN Participating Countries"
$ws.Range("B11").Value = "ປະເທດທີ່ເຂົ້າຮ່ວມ ACORN"
$ws.Range("A12").Value = "All 'orgname' are provided."
$ws.Range("B12").Value = "ທຸກໆ 'orgname'ທັງໝົດແມ່ນຖືກຕອບສະໜອງ"
$ws.Range("A13").Value = "All 'patid' are provided."
$ws.Range("B13").Value = "patid' ທັງໝົດແມ່ນຖືກຕອບສະໜອງ"
$ws.Range("A14").Value = "All 'specdate' are provided."
$ws.Range("B14").Value = "specdate' ທັງໝົດແມ່ນຖືກຕອບສະໜອງ"
$ws.Range("A15").Value = "All 'specdate' are today or before today."
$ws.Range("B15").Value = "ທຸກໆ 'specdate' ແມ່ນມື້ນີ້ ຫຼື ມື້ກ່ອນໜ້ານີ້"
$ws.Range("A16").Value = "All 'specgroup' are provided."
$ws.Range("B16").Value = "ທັງໝົດຂອງ'specgroup' ແມ່ນຖືກກຳນົດໄວ້ແລ້ວ"
$ws.Range("A17").Value = "All 'specid' are provided."
$ws.Range("B17").Value = "specid' ທັງໝົດແມ່ນຖືກຕອບສະໜອງ"
$ws.Range("A18").Value = "All dates of enrolment for HAI patients have a matching date in the HAI survey dataset"
$ws.Range("B18").Value = "ທຸກວັນທີເຂົ້າການສຶກສາສຳລັບຄົນເຈັບHAI ຈັບຄູ່ເຂົ້າກັບວັນທີເຮັດການສຳຫຼວດ HAI"
$ws.Range("A19").Value = "All Other Organisms"
$ws.Range("B19").Value = "ຕົວເຊື້ອອື່ນໆທັງໝົດ"
$ws.Range("A20").Value = "All valid records have an ACORN ID."
$ws.Range("B20").Value = "TBT"
$ws.Range("A21").Value = "AMR"
$ws.Range("B21").Value = "TBT"
$ws.Range("A22").Value = "and generate enrolment log."
$ws.Range("B22").Value = "ແລະສ້າງບັນທຶກການເຂົ້າຮ່ວມ"
$ws.Range("A23").Value = "Attempting to connect."
$ws.Range("B23").Value = "ກຳລັງພະຍາຍາມເຊື່ອມຕໍ່"
$ws.Range("A24").Value = "Blood culture collected within 24 hours of admission (CAI) / symptom onset (HAI)"
$ws.Range("B24").Value = "ປູກເລືອດພາຍໃນ24ຊົ່ວໂມງຂອງການເຂົ້ານອນ (CAI) / ເລີ້ມມີການຊຶມເຊື້ອໃນໂຮງໝໍ (HAI)"
$ws.Range("A25").Value = "Blood Culture Contaminants"
$ws.Range("B25").Value = "ການປົນເປື້ອນຂອງການປູກເລືອດ"
$ws.Range("A26").Value = "Bloodstream Infection (BSI)"
$ws.Range("B26").Value = "ການຊຶມເຊື້ອໃນກະແສເລືອດ"
$ws.Range("A27").Value = "Calculated age is consistent with 'Age Category'"
$ws.Range("B27").Value = "ການຄິດໄລ່ອາຍຸແມ່ນເໝາະສົມກັບ 'Age Category'"
$ws.Range("A28").Value = "Calculated age isn't always consistent with 'Age Category'"
$ws.Range("B28").Value = "ການຄິດໄລ່ອາຍຸອາດບໍ່ເໝາະສົມກັບ 'Age Category' ສະເໝີໄປ"
$ws.Range("A29").Value = "Cancel"
$ws.Range("B29").Value = "ຍົກເລີກ"
$ws.Range("A30").Value = "Care should be taken when interpreting rates and AMR profiles where there are small numbers of cases or bacterial isolates: point estimates may be unreliable."
$ws.Range("B30").Value = "ຄວນໃຊ້ຄວາມລະມັດລະວັງໃນການແປຄວາມໝາຍອັດຕາແລະຮູບຮ່າງລັກສະນະຂອງການຕ້ານຕໍ່ຢາຕ້ານເຊື້ອທີ່ມີຈຳນວນຄົນເຈັບໜ້ອຍ ຫຼື ແບັກທີເຣຍທີ່ແຍກໄດ້: ການແປຜົນອາດບໍ່ໜ້າເຊື່ອຖື"
$ws.Range("A31").Value = "Clinical and day-28 outcomes are consistent."
$ws.Range("B31").Value = "ອາການຄຣີນິກ ແລະ ການຕິດຕາມມື້ທີ 28 ແມ່ນກົງກັນ"
$ws.Range("A32").Value = "Clinical and day-28 outcomes aren't consistent for some dead patients."
$ws.Range("B32").Value = "ອາການຄຣີນິກ ແລະ ການຕິດຕາມມື້ທີ 28 ບໍ່ກົງກັນໃນບາງຄົນເຈັບທີ່ເສຍຊີວິດ"
$ws.Range("A33").Value = "Clinical Outcome"
$ws.Range("B33").Value = "ອາການຫຼັງຕິດຕາມ"
$ws.Range("A34").Value = "Clinical Outcome Status:"
$ws.Range("B34").Value = "ສະຖານະພາບອາການຄຣີນິກຫຼັງຕິດຕາມ"
$ws.Range("A35").Value = "Co-resistances"
$ws.Range("B35").Value = "TBT"
$ws.Range("A36").Value = "Combine Susceptible + Intermediate"
$ws.Range("B36").Value = "Susceptible + Intermediate ລວມເຂົ້າກັນ"
$ws.Range("A37").Value = "Consider saving .acorn file on the cloud for additional security."
$ws.Range("B37").Value = "ຄວນພິຈາລະນາບັນທຶກຂໍ້ມູນ.acorn ໃນ cloud ເພື່ອຄວາມປອດໄພເພີ້ມເຕີມ"
$ws.Range("A38").Value = "Contains names of organisms before and after mapping."
$ws.Range("B38").Value = "TBT"
$ws.Range("A39").Value = "Couldn't connect to server. Please check internet access."
$ws.Range("B39").Value = "ບໍ່ສາມາດເຊື່ອມຕໍ່ກັບ server ໄດ້. ກະລຸນາກວດສອບການເຂົ້າເຖິງອີນເຕີເນັດ"
$ws.Range("A40").Value = "Critical errors with clinical data."
$ws.Range("B40").Value = "ຂໍ້ຜິດພາດທີ່ຮ້າຍແຮງກ່ຽວກັບຂໍ້ມູນທາງຄຣີນິກ."
$ws.Range("A41").Value = "Culture results per specimen type"
$ws.Range("B41").Value = "ຜົນການປູກຕໍ່ກັບປະເພດຕົວຢ່າງ"
$ws.Range("A42").Value = "Data Management"
$ws.Range("B42").Value = "ການຈັດການຂໍ້ມູນ"
$ws.Range("A43").Value = "Date of Enrolment"
$ws.Range("B43").Value = "ວັນທີເຂົ້າຮ່ວມການສຶກສາ"
$ws.Range("A44").Value = "Day 28"
$ws.Range("B44").Value = "ມື້ທີ28"
$ws.Range("A45").Value = "Day 28 Status:"
$ws.Range("B45").Value = "ສະຖານະພາບມື້ທີ28"
$ws.Range("A46").Value = "Diagnosis at Enrolment"
$ws.Range("B46").Value = "ການບົ່ງມະຕິເວລາເຂົ້າການສຶກສາ"
$ws.Range("A47").Value = "Dismiss"
$ws.Range("B47").Value = "TBT"
$ws.Range("A48").Value = "Distribution of Enrolments"
$ws.Range("B48").Value = "ການແຈກຢາຍຂອງການເຂົ້າຮ່ວມ"
$ws.Range("A49").Value = "Download Enrolment Log (.xlsx)"
$ws.Range("B49").Value = "ດຶງຂໍ້ມູນບັນທຶກການເຂົ້າຮ່ວມ (.xlsx)"
$ws.Range("A50").Value = "Download Lab Log (.xlsx)"
$ws.Range("B50").Value = "TBT"
$ws.Range("A70").Value = "HAI point prevalence by "
$ws.Range("B70").Value = "TBT"
$ws.Range("A110").Value = "Remove 'Not Cultured' specimens"
$ws.Range("B110").Value = "TBT"
$ws.Range("A111").Value = "Remove blood culture contaminants from the following visualizations"
$ws.Range("B111").Value = "ການປູກເລືອດທີ່ມີການປົນເປື້ອນແມ່ນລົບຜົນອອກຈາກຜົນການສະແດງຂໍ້ມູນ"
$ws.Range("A112").Value = "Reset Enrolments Filters"
$ws.Range("B112").Value = "Reset Enrolments Filters"
$ws.Range("A113").Value = "Resistance to 3rd gen. Cephalosporins Over Time"
$ws.Range("B113").Value = "ການຕ້ານຕໍ່ 3rd gen. Cephalosporins ຄ່ອຍເປັນຄ່ອຍໄປຕາມໄລຍະເວລາ"
$ws.Range("A114").Value = "Resistance to Carbapenems Over Time"
$ws.Range("B114").Value = "ມີການຕ້ານຕໍ່ຢາ Carbapenems ຄ່ອຍເປັນຄ່ອຍໄປຕາມໄລຍະເວລາ"
$ws.Range("A115").Value = "Resistance to Fluoroquinolones Over Time"
$ws.Range("B115").Value = "ການຕ້ານຕໍ່ Fluoroquinolones ຄ່ອຍເປັນຄ່ອຍໄປຕາມໄລຍະເວລາ"
$ws.Range("A116").Value = "Resistance to Oxacillin Over Time"
$ws.Range("B116").Value = "ການຕ້ານຕໍ່ຢາ Oxacillin ເປັນໄປຕາມໄລຍະເວລາ"
$ws.Range("A117").Value = "Resistance to Penicillin G - meningitis Over Time"
$ws.Range("B117").Value = "ການຕ້ານຕໍ່ຢາ Penicillin G- meningitis ເປັນໄປຕາມໄລຍະເວລາ"
$ws.Range("A118").Value = "Resistance to Penicillin G Over Time"
$ws.Range("B118").Value = "ການຕ້ານຕໍ່ຢາ Penicillin G ເປັນໄປຕາມໄລຍະເວລາ"
$ws.Range("A119").Value = "Retriving data from REDCap server."
$ws.Range("B119").Value = "ກຳລັງກູ້ຄືນຂໍ້ມູນຈາກຖານ REDCap"
$ws.Range("A120").Value = "Save .acorn file"
$ws.Range("B120").Value = "ບັນທຶກເອກະສານຂອງ.acorn"
$ws.Range("A121").Value = "Save acorn data"
$ws.Range("B121").Value = "ບັນທຶກຂໍ້ມູນ acorn"
$ws.Range("A122").Value = "Save on Server"
$ws.Range("B122").Value = "ບັນທຶກລົງໃນຖານຂໍ້ມູນ"
$ws.Range("A123").Value = "See Breakdown by Ward"
$ws.Range("B123").Value = "ເບິ່ງລາຍລະອຽດໂດຍອິງໃສ່ພະແນກ"
$ws.Range("A124").Value = "See by Week"
$ws.Range("B124").Value = "ເບິງລາຍອາທິດ"
